# Fix the URL to the Eresume:
#  - drop the old "_GoBack" bookmark near the top of the document
#  - accept the tracked deletion of "com/MrChido" (the old, broken path)
#    so only the two insertions ("mrchido." and "io/") remain
#  - re-create the "_GoBack" bookmark at the end of the "SKILLS" paragraph,
#    which is where Word left it after the most recent edit

$d = $word.ActiveDocument

# 1. Remove the stale "_GoBack" bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Accept the tracked deletion revision that removed "com/MrChido" from
#    the Eresume hyperlink, leaving the corrected "mrchido.github.io/" URL.
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $rev = $d.Revisions.Item($i)
    if ($rev.Type -eq 2) {
        $rev.Accept()
    }
}

# 3. Re-add "_GoBack" right after the "SKILLS" heading text, at the end of
#    that paragraph (its new location after the edit above).
$skills = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "SKILLS") {
        $skills = $p
        break
    }
}

$endRng = $skills.Range
$endRng.Collapse(0)
$endRng.Text = "X"

$markerStart = $skills.Range.End - 2
$mark = $d.Range($markerStart, $markerStart)
$d.Bookmarks.Add("_GoBack", $mark)

$placeholder = $d.Range($markerStart, $markerStart + 1)
$placeholder.Text = ""
